$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: informational note about screws needed when using the case.
# Formatting: same font as the header/quantity cells (Arial, theme color),
# but with no border - achieved by cloning A1's format then removing the border.
$ws.Range("A7").Value = "If using the case you need these screws"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Borders.LineStyle = -4142

# Row 8: the screw part itself, styled like the other hyperlinked part names
# (A2:A5), plus an empty bordered cell in column B matching the Quantity column.
$ws.Range("A8").Value = "M4x25MM "

# Hyperlink for the new screw line item (added before formatting, since
# creating a hyperlink re-styles the cell with the built-in Hyperlink style).
$ws.Hyperlinks.Add($ws.Range("A8"), "https://www.mcmaster.com/91290A247/") | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").ClearContents()
